$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-deleted last data row (old row 10)
$ws.Rows(10).Delete()

# Overwrite data rows 2-9 with the recomputed TPM values
# Row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Hp'
$ws.Range("C2").Value = 'Cd163'
$ws.Range("D2").Value = 'MuSCs'
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1238986666666667
$ws.Range("H2").Value = 0.371696
$ws.Range("I2").Value = 0.01923905185495286
$ws.Range("J2").Value = 0.01923905185495286
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.690195
$ws.Range("N2").Value = 2.070585
$ws.Range("O2").Value = 0.1060250152438306
$ws.Range("P2").Value = 0.1060250152438306
$ws.Range("Q2").Value = 0.08551424024
$ws.Range("R2").Value = 0.76962816216
$ws.Range("S2").Value = 0.002039820766198225
$ws.Range("T2").Value = 0.002039820766198225

# Row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Hp'
$ws.Range("C3").Value = 'Cd163'
$ws.Range("D3").Value = 'Resolving-Mac'
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1238986666666667
$ws.Range("H3").Value = 0.371696
$ws.Range("I3").Value = 0.01923905185495286
$ws.Range("J3").Value = 0.01923905185495286
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.819542333333334
$ws.Range("N3").Value = 17.458627
$ws.Range("O3").Value = 0.8939749847561693
$ws.Range("P3").Value = 0.8939749847561693
$ws.Range("Q3").Value = 0.7210335357102223
$ws.Range("R3").Value = 6.489301821392001
$ws.Range("S3").Value = 0.01719923108875463
$ws.Range("T3").Value = 0.01719923108875463

# Row 4
$ws.Range("A4").Value = 'FAPs'
$ws.Range("B4").Value = 'Hp'
$ws.Range("C4").Value = 'Cd163'
$ws.Range("D4").Value = 'MuSCs'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.810518333333333
$ws.Range("H4").Value = 17.431555
$ws.Range("I4").Value = 0.9022604239955847
$ws.Range("J4").Value = 0.9022604239955845
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.690195
$ws.Range("N4").Value = 2.070585
$ws.Range("O4").Value = 0.1060250152438306
$ws.Range("P4").Value = 0.1060250152438306
$ws.Range("Q4").Value = 4.010390701075
$ws.Range("R4").Value = 36.093516309675
$ws.Range("S4").Value = 0.09566217520803695
$ws.Range("T4").Value = 0.09566217520803691

# Row 5
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Hp'
$ws.Range("C5").Value = 'Cd163'
$ws.Range("D5").Value = 'Resolving-Mac'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.810518333333333
$ws.Range("H5").Value = 17.431555
$ws.Range("I5").Value = 0.9022604239955847
$ws.Range("J5").Value = 0.9022604239955845
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.819542333333334
$ws.Range("N5").Value = 17.458627
$ws.Range("O5").Value = 0.8939749847561693
$ws.Range("P5").Value = 0.8939749847561693
$ws.Range("Q5").Value = 33.81455741944278
$ws.Range("R5").Value = 304.331016774985
$ws.Range("S5").Value = 0.8065982487875477
$ws.Range("T5").Value = 0.8065982487875475

# Row 6
$ws.Range("A6").Value = 'MuSCs'
$ws.Range("B6").Value = 'Hp'
$ws.Range("C6").Value = 'Cd163'
$ws.Range("D6").Value = 'MuSCs'
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1062546666666667
$ws.Range("H6").Value = 0.318764
$ws.Range("I6").Value = 0.01649928200866351
$ws.Range("J6").Value = 0.01649928200866351
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.690195
$ws.Range("N6").Value = 2.070585
$ws.Range("O6").Value = 0.1060250152438306
$ws.Range("P6").Value = 0.1060250152438306
$ws.Range("Q6").Value = 0.07333643966
$ws.Range("R6").Value = 0.6600279569399999
$ws.Range("S6").Value = 0.001749336626480809
$ws.Range("T6").Value = 0.001749336626480809

# Row 7
$ws.Range("A7").Value = 'MuSCs'
$ws.Range("B7").Value = 'Hp'
$ws.Range("C7").Value = 'Cd163'
$ws.Range("D7").Value = 'Resolving-Mac'
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1062546666666667
$ws.Range("H7").Value = 0.318764
$ws.Range("I7").Value = 0.01649928200866351
$ws.Range("J7").Value = 0.01649928200866351
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.819542333333334
$ws.Range("N7").Value = 17.458627
$ws.Range("O7").Value = 0.8939749847561693
$ws.Range("P7").Value = 0.8939749847561693
$ws.Range("Q7").Value = 0.6183535307808889
$ws.Range("R7").Value = 5.565181777028
$ws.Range("S7").Value = 0.0147499453821827
$ws.Range("T7").Value = 0.0147499453821827

# Row 8
$ws.Range("A8").Value = 'Resolving-Mac'
$ws.Range("B8").Value = 'Hp'
$ws.Range("C8").Value = 'Cd163'
$ws.Range("D8").Value = 'MuSCs'
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3992853333333333
$ws.Range("H8").Value = 1.197856
$ws.Range("I8").Value = 0.06200124214079897
$ws.Range("J8").Value = 0.06200124214079897
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.690195
$ws.Range("N8").Value = 2.070585
$ws.Range("O8").Value = 0.1060250152438306
$ws.Range("P8").Value = 0.1060250152438306
$ws.Range("Q8").Value = 0.27558474064
$ws.Range("R8").Value = 2.48026266576
$ws.Range("S8").Value = 0.006573682643114645
$ws.Range("T8").Value = 0.006573682643114643

# Row 9
$ws.Range("A9").Value = 'Resolving-Mac'
$ws.Range("B9").Value = 'Hp'
$ws.Range("C9").Value = 'Cd163'
$ws.Range("D9").Value = 'Resolving-Mac'
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3992853333333333
$ws.Range("H9").Value = 1.197856
$ws.Range("I9").Value = 0.06200124214079897
$ws.Range("J9").Value = 0.06200124214079897
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.819542333333334
$ws.Range("N9").Value = 17.458627
$ws.Range("O9").Value = 0.8939749847561693
$ws.Range("P9").Value = 0.8939749847561693
$ws.Range("Q9").Value = 2.323657900412444
$ws.Range("R9").Value = 20.912921103712
$ws.Range("S9").Value = 0.05542755949768433
$ws.Range("T9").Value = 0.05542755949768432

